# Remove the slide titled "Набор статусов и переходов жизненного цикла задачи"
# (sldId 282 in the original deck). All other slides/content are left untouched.
$p = $ppt.ActivePresentation

$targetTitle = "Набор статусов и переходов жизненного цикла задачи"
$targetIndex = -1

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.Count -gt 0) {
        $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
        if ($title -eq $targetTitle) {
            $targetIndex = $i
        }
    }
}

if ($targetIndex -gt 0) {
    $p.Slides.Item($targetIndex).Delete()
}
